# Update the "想去人数" (wanted-to-go count) figures in column F across
# the "展览", "演出" and "全部类型" sheets, matching refreshed site data.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 31
$ws1.Range("F3").Value = 8971
$ws1.Range("F4").Value = 2702
$ws1.Range("F5").Value = 959
$ws1.Range("F8").Value = 700
$ws1.Range("F9").Value = 128
$ws1.Range("F10").Value = 78
$ws1.Range("F13").Value = 3834
$ws1.Range("F14").Value = 290
$ws1.Range("F15").Value = 172
$ws1.Range("F22").Value = 1389
$ws1.Range("F24").Value = 481
$ws1.Range("F27").Value = 174
$ws1.Range("F28").Value = 369
$ws1.Range("F29").Value = 70
$ws1.Range("F33").Value = 718
$ws1.Range("F34").Value = 56
$ws1.Range("F37").Value = 4
$ws1.Range("F39").Value = 3
$ws1.Range("F41").Value = 197
$ws1.Range("F42").Value = 342
$ws1.Range("F43").Value = 25
$ws1.Range("F44").Value = 16

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 3

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 31
$ws4.Range("F6").Value = 700
$ws4.Range("F7").Value = 128
$ws4.Range("F8").Value = 78
$ws4.Range("F12").Value = 3834
$ws4.Range("F13").Value = 290
$ws4.Range("F14").Value = 172
$ws4.Range("F16").Value = 3
$ws4.Range("F17").Value = 803
$ws4.Range("F26").Value = 1389
$ws4.Range("F28").Value = 481
$ws4.Range("F31").Value = 174
$ws4.Range("F33").Value = 369
$ws4.Range("F34").Value = 70
$ws4.Range("F37").Value = 718
$ws4.Range("F38").Value = 56
$ws4.Range("F41").Value = 4
$ws4.Range("F43").Value = 3
$ws4.Range("F44").Value = 197
$ws4.Range("F45").Value = 342
$ws4.Range("F46").Value = 25
$ws4.Range("F47").Value = 16
